$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the title in A1
$ws.Range("A1").Value = "Going Up!"

# Fill in the previously-empty "Value" column (D) entries
$ws.Range("D5").Value = 15
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 10
$ws.Range("D8").Value = 5
$ws.Range("D11").Value = 5
$ws.Range("D13").Value = 10

# Update the active selection to match the target view
$ws.Range("D12").Select() | Out-Null
